$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely (no naive-forecast value available for these rows)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Update values that changed due to the naive component forecaster bug fix
$ws.Range("E3").Value  = -5.168396053267521
$ws.Range("C4").Value  = -3.956152295564885
$ws.Range("C5").Value  = 1.234995474941436
$ws.Range("C6").Value  = 0.8993608108207818
$ws.Range("E6").Value  = 0.4755443417510108
$ws.Range("C8").Value  = 0.02019328874802717
$ws.Range("E8").Value  = -1.194807813319176
$ws.Range("E10").Value = -1.02250637024307
$ws.Range("E11").Value = -0.7240982069265045
$ws.Range("C12").Value = 0.0720185131838802
$ws.Range("E12").Value = 1.255028673974068
$ws.Range("E13").Value = -3.305525567352907
$ws.Range("C14").Value = -0.8017595264762423
$ws.Range("E15").Value = 11.45073880931156
$ws.Range("C16").Value = 0.9704846793491706
$ws.Range("E16").Value = -1.172596637408219
$ws.Range("C18").Value = 0.3928252664241683
$ws.Range("C19").Value = 0.3224026462283369
$ws.Range("E19").Value = -3.037731958703715
